# Weekly fruit/vegetable price update:
# Insert a new price record as row 74 (shifting the existing rows 74-79 down
# to 75-80), matching the new "Argentina(o)" Alcachofa entry for
# Vega Monumental Concepción.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 74:79 down to 75:80 and leave a blank row 74 for the new record.
$ws.Rows.Item(74).Insert()

# Populate the newly inserted row 74 with the new record's data.
$ws.Cells.Item(74, 1).Value = 11
$ws.Cells.Item(74, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(74, 3).Value = "Bíobío"
$ws.Cells.Item(74, 4).Value = 45021
$ws.Cells.Item(74, 5).Value = 8
$ws.Cells.Item(74, 6).Value = 100112013
$ws.Cells.Item(74, 7).Value = "Alcachofa"
$ws.Cells.Item(74, 8).Value = "Argentina(o)"
$ws.Cells.Item(74, 9).Value = "Primera"
$ws.Cells.Item(74, 10).Value = 50
$ws.Cells.Item(74, 11).Value = 16000
$ws.Cells.Item(74, 12).Value = 17000
$ws.Cells.Item(74, 13).Value = 16600
$ws.Cells.Item(74, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(74, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(74, 16).Value = 332
$ws.Cells.Item(74, 17).Value = 50
$ws.Cells.Item(74, 18).Value = "Hortaliza"
